# Applies the "Correzione problema su update category admin" edit described
# by the diff:
#   1. Rewords / restructures the "Lavori minoritari" block under
#      "Lista lavori Giacomo" (new text for "Finire pdf/fattura" and
#      "Sistemare menù categorie", reorders "Termina traduzioni alt HTML" /
#      "Trova privacy in inglese" below the other bullet points, moves the
#      page-break + bookmark + "Lista lavori Francesco:" heading into its own
#      paragraph).
#   2. Adds a new "Menù categorie." bullet right after "Lavori con massima
#      priorità:" (numId 3).
#   3. Removes the two questions ("Chiedere livelli di sotto categorie..."
#      and "Vedere anche pdf fattura/dettaglio(??).") from the "Domande per
#      cliente:" block.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, [string]$text) {
    $target = $text.Trim()
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t.Trim() -eq $target) {
            return $i
        }
    }
    return -1
}

function Get-ParaIndexByPrefix($doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function New-OpenXmlPackage([string]$body) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$body
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------------
# Edit 1: "Finire pdf/fattura." ... "Non fa comparire...Lista lavori
# Francesco:" -- rewrite the whole run of bullets in one shot so the
# reordering / restructuring matches exactly.
# ---------------------------------------------------------------------------

$startIdx = Get-ParaIndexByText $d "Finire pdf/fattura."
$endIdx = Get-ParaIndexByPrefix $d "Non fa comparire sempre i puntini nelle descrizioni, ma solo quando superano X caratteri."

if ($startIdx -eq -1 -or $endIdx -eq -1) {
    throw "Could not locate the 'Lavori minoritari' block anchors ($startIdx, $endIdx)"
}

$pStart = $d.Paragraphs.Item($startIdx)
$pEnd = $d.Paragraphs.Item($endIdx)
$blockRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)

$newBlockBody = @'
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="4BACC6" w:themeColor="accent5"/></w:rPr></w:pPr><w:r><w:t>Finire pdf/fattura</w:t></w:r><w:r><w:t xml:space="preserve"> (vuole anche file solo per dettaglio, da dinamicizzare il contesto)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="4BACC6" w:themeColor="accent5"/></w:rPr></w:pPr><w:r><w:t>Sistemare menù categorie: 3 vanno bene, vediamo come implementare la cosa.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="4BACC6" w:themeColor="accent5"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Mettere pager in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/site/list*</w:t></w:r></w:p>
<w:p/>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Lavori minoritari:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Barra di ricerca sui </w:t></w:r><w:r><w:t>prodotti</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>precedenza pari a 0)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Migliorare </w:t></w:r><w:r><w:t>usabilità</w:t></w:r><w:r><w:t xml:space="preserve"> sezione company.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Non fa comparire sempre i puntini nelle descrizioni, ma solo quando superano X caratteri.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Termina traduzioni alt HTML</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Trova privacy in inglese</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:ind w:left="1080"/></w:pPr><w:r><w:br w:type="page"/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:color w:val="F79646" w:themeColor="accent6"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Lista lavori Francesco:</w:t></w:r></w:p>
'@

$blockRange.InsertXML((New-OpenXmlPackage $newBlockBody)) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: insert "Menù categorie." right after "Lavori con massima
# priorità:" (the numId 3 list, the second occurrence of that heading).
# ---------------------------------------------------------------------------

$priorityIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]7).Trim()
    if ($t -eq "Lavori con massima priorità:") {
        # Disambiguate from the similarly-worded "Lista lavori Giacomo"
        # heading near the top of the document (which carries extra runs on
        # the same paragraph and is followed by "Implementazione
        # MultiLanguage..."): the numId-3 heading targeted by the diff is
        # immediately followed by "Lavori minoritari ma prioritari:".
        $nextText = ""
        if ($i -lt $d.Paragraphs.Count) {
            $nextText = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd([char]13, [char]7).Trim()
        }
        if ($nextText -eq "Lavori minoritari ma prioritari:") {
            $priorityIdx = $i
            break
        }
    }
}

if ($priorityIdx -eq -1) {
    throw "Could not locate the 'Lavori con massima priorità:' (numId 3) heading"
}

$pPriority = $d.Paragraphs.Item($priorityIdx)
$priorityRange = $d.Range($pPriority.Range.Start, $pPriority.Range.End)

$priorityBody = @'
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Lavori con massima priorità:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Menù categorie.</w:t></w:r></w:p>
'@

$priorityRange.InsertXML((New-OpenXmlPackage $priorityBody)) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: remove the two "Domande per cliente:" sub-bullets about category
# depth and the pdf invoice detail.
# ---------------------------------------------------------------------------

$q1Idx = Get-ParaIndexByText $d "Chiedere livelli di sotto categorie (per farsi un idea), possibile soluzione: mostro profondità di due e poi nella pagina faccio vedere tutte le sotto categorie."
$q2Idx = Get-ParaIndexByText $d "Vedere anche pdf fattura/dettaglio(??)."

if ($q1Idx -eq -1 -or $q2Idx -eq -1) {
    throw "Could not locate the 'Domande per cliente' questions to remove ($q1Idx, $q2Idx)"
}
if ($q2Idx -ne $q1Idx + 1) {
    throw "Unexpected paragraph order around 'Domande per cliente' ($q1Idx, $q2Idx)"
}

$pQ1 = $d.Paragraphs.Item($q1Idx)
$pQ2 = $d.Paragraphs.Item($q2Idx)
$qRange = $d.Range($pQ1.Range.Start, $pQ2.Range.End)
$qRange.Delete() | Out-Null
